# Update the students workbook:
#  - A1 used to hold a leftover "students.xlsx" label (small green Consolas
#    comment-style font); replace it with the real "Student ID" column
#    header and give it the same look as the other header cells (B1/C1).
#  - Leave the active selection on E4, matching the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Student ID"

# Match A1's formatting to the existing header style used by B1/C1
# (bold Calibri, thin border, centered) instead of leaving the old
# ad-hoc "comment" font on it.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Restore/save the active cell selection.
$ws.Range("E4").Select() | Out-Null
